$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "CDB"

$ws.Range("T1").Value = 15010
$ws.Range("U1").Value = 503

$ws.Range("E9").Value = 12301
$ws.Range("D9").Value = $null
$ws.Range("G9").Value = 12026
$ws.Range("F9").Value = $null
$ws.Range("H9").Value = 12025
$ws.Range("I9").Value = 12024
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 402
$ws.Range("L9").Value = 404
$ws.Range("M9").Value = 406
$ws.Range("N9").Value = 408
$ws.Range("O9").Value = 410
$ws.Range("P9").Value = 412
$ws.Range("Q9").Value = 414
$ws.Range("R9").Value = 416
$ws.Range("S9").Value = 418
$ws.Range("T9").Value = 420
$ws.Range("U9").Value = 422
$ws.Range("V9").Value = 424
$ws.Range("W9").Value = 426
$ws.Range("X9").Value = 428

$ws.Range("D10").Value = 12016
$ws.Range("E10").Value = 12302
$ws.Range("F10").Value = 12023
$ws.Range("G10").Value = 12022
$ws.Range("H10").Value = 12021
$ws.Range("I10").Value = 12020
$ws.Range("J10").Value = 401
$ws.Range("K10").Value = 403
$ws.Range("L10").Value = 405
$ws.Range("M10").Value = 407
$ws.Range("N10").Value = 409
$ws.Range("O10").Value = 411
$ws.Range("P10").Value = 413
$ws.Range("Q10").Value = 415
$ws.Range("R10").Value = 417
$ws.Range("S10").Value = 419
$ws.Range("T10").Value = 421
$ws.Range("U10").Value = 423
$ws.Range("V10").Value = 425
$ws.Range("W10").Value = 427
$ws.Range("X10").Value = 429

$ws.Range("E11").Value = 12303
$ws.Range("D11").Value = $null

$ws.Range("E16").Clear()
$ws.Range("G16").Clear()
$ws.Range("I16").ClearFormats()

$ws.Range("T1").Select()
